# Auto-generated cell updates based on the provided diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Force the cell to hold the literal text given, even if it
    # looks like a number, and leave no residual explicit style
    # behind (matches original inlineStr cells with no s= attribute).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '63.543.85'
Set-TextValue $ws.Range('E2') '  +0.14%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.616.61'
Set-TextValue $ws.Range('E3') '  -0.60%  '

# Row 4
Set-TextValue $ws.Range('E4') '  -0.09%  '

# Row 5
Set-TextValue $ws.Range('D5') '594.98'
Set-TextValue $ws.Range('E5') '  -1.82%  '

# Row 6
Set-TextValue $ws.Range('D6') '150.15'
Set-TextValue $ws.Range('E6') '  +2.01%  '

# Row 7
Set-TextValue $ws.Range('E7') '  -0.05%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.589'
Set-TextValue $ws.Range('E8') '  -0.02%  '

# Row 9
Set-TextValue $ws.Range('E9') '  +0.08%  '

# Row 10
Set-TextValue $ws.Range('D10') '5.68'
Set-TextValue $ws.Range('E10') '  +1.87%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.382'
Set-TextValue $ws.Range('E11') '  +2.75%  '

# Row 12
Set-TextValue $ws.Range('E12') '  -1.12%  '

# Row 13
Set-TextValue $ws.Range('D13') '27.65'
Set-TextValue $ws.Range('E13') '  +0.40%  '

# Row 14
Set-TextValue $ws.Range('D14') '3.086.23'
Set-TextValue $ws.Range('E14') '  -0.75%  '

# Row 15
Set-TextValue $ws.Range('D15') '63.367.28'
Set-TextValue $ws.Range('E15') '  +0.10%  '

# Row 16
Set-TextValue $ws.Range('E16') '  +2.38%  '

# Row 17
Set-TextValue $ws.Range('D17') '2.619.69'
Set-TextValue $ws.Range('E17') '  -1.39%  '

# Row 18
Set-TextValue $ws.Range('D18') '12.34'
Set-TextValue $ws.Range('E18') '  +6.76%  '

# Row 19
Set-TextValue $ws.Range('E19') '  +1.99%  '

# Row 20
Set-TextValue $ws.Range('D20') '346.38'
Set-TextValue $ws.Range('E20') '  +0.44%  '

# Row 21
Set-TextValue $ws.Range('E21') '  -0.49%  '

# Row 22
Set-TextValue $ws.Range('D22') '0.998'
Set-TextValue $ws.Range('E22') '  -0.23%  '

# Row 23
Set-TextValue $ws.Range('E23') '  +2.99%  '

# Row 24
Set-TextValue $ws.Range('D24') '66.49'
Set-TextValue $ws.Range('E24') '  -0.67%  '

# Row 25
Set-TextValue $ws.Range('E25') '  +9.47%  '

# Row 26
Set-TextValue $ws.Range('D26') '9.18'
Set-TextValue $ws.Range('E26') '  +1.26%  '

# Row 27
Set-TextValue $ws.Range('E27') '  -1.79%  '

# Row 28
Set-TextValue $ws.Range('B28') 'Bittensor'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D28') '552.89'
Set-TextValue $ws.Range('E28') '  -2.10%  '

# Row 29
Set-TextValue $ws.Range('B29') 'Aptos'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D29') '8.14'
Set-TextValue $ws.Range('E29') '  +1.10%  '

# Row 30
Set-TextValue $ws.Range('D30') '0.162'
Set-TextValue $ws.Range('E30') '  -0.93%  '

# Row 31
Set-TextValue $ws.Range('D31') '0.998'
Set-TextValue $ws.Range('E31') '  -0.19%  '

# Row 32
Set-TextValue $ws.Range('D32') '2.04'
Set-TextValue $ws.Range('E32') '  -0.91%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.0₃0845'
Set-TextValue $ws.Range('E33') '  -0.79%  '

# Row 34
Set-TextValue $ws.Range('E34') '  -0.74%  '

# Row 35
Set-TextValue $ws.Range('E35') '  +1.07%  '

# Row 36
Set-TextValue $ws.Range('D36') '168.07'
Set-TextValue $ws.Range('E36') '  +0.36%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.414'
Set-TextValue $ws.Range('E37') '  +1.94%  '

# Row 38
Set-TextValue $ws.Range('E38') '  -0.03%  '

# Row 39
Set-TextValue $ws.Range('D39') '19.44'
Set-TextValue $ws.Range('E39') '  +1.52%  '

# Row 40
Set-TextValue $ws.Range('E40') '  -1.54%  '

# Row 41
Set-TextValue $ws.Range('D41') '1.00'
Set-TextValue $ws.Range('E41') '  +0.03%  '

# Row 42
Set-TextValue $ws.Range('D42') '166.52'
Set-TextValue $ws.Range('E42') '  +0.31%  '

# Row 43
Set-TextValue $ws.Range('D43') '39.72'
Set-TextValue $ws.Range('E43') '  -1.02%  '

# Row 44
Set-TextValue $ws.Range('D44') '3.91'
Set-TextValue $ws.Range('E44') '  +3.11%  '

# Row 45
Set-TextValue $ws.Range('E45') '  +2.68%  '

# Row 46
Set-TextValue $ws.Range('D46') '21.52'
Set-TextValue $ws.Range('E46') '  -2.91%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.630'
Set-TextValue $ws.Range('E47') '  -0.01%  '

# Row 48
Set-TextValue $ws.Range('E48') '  +0.88%  '

# Row 49
Set-TextValue $ws.Range('B49') 'BabyDogeCoin'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D49') '0.0₆0252'
Set-TextValue $ws.Range('E49') '  +26.25%  '

# Row 50
Set-TextValue $ws.Range('B50') 'dogwifhat'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D50') '1.99'
Set-TextValue $ws.Range('E50') '  +3.18%  '

# Row 51
Set-TextValue $ws.Range('B51') 'Stellar'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D51') '0.0966'
Set-TextValue $ws.Range('E51') '  +0.43%  '
